# Apply the "35Ca" sheet update: add new analysis columns X (measured value),
# Y (uncertainty / count) and Z (= X*100/100.6, a 100/100.6 normalisation)
# to the existing nuclear-data table on sheet "35Ca".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("35Ca")

# ---------------------------------------------------------------------------
# Number formats used by the new columns:
#   fmt1 ("0.0")  -> mirrors the workbook's existing "0.0" style, applied to
#                    the "larger" values
#   fmt2 ("0.00") -> a new, more precise format used for the small values
# Column Y (counts) keeps the plain/general format already inherited from the
# sheet's default column style.
# ---------------------------------------------------------------------------
$fmt1 = "0.0"
$fmt2 = "0.00"

function SetXYZ {
    param($row, $xval, $yval, $fmt)
    $ws.Range("X$row").Value = $xval
    $ws.Range("X$row").NumberFormat = $fmt
    $ws.Range("Y$row").Value = $yval
    $ws.Range("Z$row").Formula = "=X$row*100/100.6"
    $ws.Range("Z$row").NumberFormat = $fmt
}

function SetEmptyXZ {
    param($row, $fmt)
    $ws.Range("X$row").NumberFormat = $fmt
    $ws.Range("Z$row").NumberFormat = $fmt1
}

# --- First data block (rows 6-19) -------------------------------------------
SetXYZ 6  48.5 1.3 $fmt1
SetXYZ 7  6    5   $fmt1
SetXYZ 8  3    3   $fmt1
SetXYZ 9  3.8  3   $fmt1
SetXYZ 10 2.9  3   $fmt1
SetXYZ 11 2.9  3   $fmt1
SetXYZ 12 4.2  4   $fmt1
SetXYZ 13 3.9  3   $fmt1
SetXYZ 14 0.72 18  $fmt2
SetXYZ 15 0.61 15  $fmt2
SetXYZ 16 1.43 17  $fmt2
SetXYZ 17 1.4  19  $fmt2
SetXYZ 18 3.8  2   $fmt1
SetXYZ 19 0.41 6   $fmt2

# Row 20 is a blank spacer row - only formatted placeholders remain.
SetEmptyXZ 20 $fmt2

# --- Second data block (rows 21-33) -----------------------------------------
SetXYZ 21 2.2  3  $fmt2
SetXYZ 22 1.09 17 $fmt2
SetXYZ 23 1.1  2  $fmt2
SetXYZ 24 2.2  3  $fmt2
SetXYZ 25 1.09 17 $fmt2
SetXYZ 26 1.1  2  $fmt2
SetXYZ 27 2.2  3  $fmt2
SetXYZ 28 1.09 17 $fmt2
SetXYZ 29 1.1  2  $fmt2

# Row 30 is a blank spacer row.
SetEmptyXZ 30 $fmt2

SetXYZ 31 8.4 6 $fmt2
SetXYZ 32 8.4 6 $fmt2
SetXYZ 33 8.4 6 $fmt2

# New blank rows 34-35 inserted after row 33 (formatted placeholders only).
SetEmptyXZ 34 $fmt2
SetEmptyXZ 35 $fmt2

# Rows 36-38 (header / reference rows) get formatted placeholders too.
SetEmptyXZ 36 $fmt2
SetEmptyXZ 37 $fmt2
SetEmptyXZ 38 $fmt2

# --- Row 39 ------------------------------------------------------------
SetXYZ 39 4.2 3 $fmt2

# ---------------------------------------------------------------------------
# Final selection, matching the author's last cursor position.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("Z6").Select()
